$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "Förändrad" (column C) date value from 2023-10-22 (45221)
# to 2023-10-25 (45224) for rows 2 through 8.
for ($row = 2; $row -le 8; $row++) {
    $ws.Cells.Item($row, 3).Value = 45224
}
